$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 8428
$ws.Range("F11").Value = 7719
$ws.Range("F21").Value = 283
$ws.Range("F25").Value = 780
$ws.Range("F28").Value = 1344
$ws.Range("F29").Value = 531
$ws.Range("F30").Value = 473
$ws.Range("F32").Value = 254
$ws.Range("F33").Value = 35
$ws.Range("F34").Value = 73
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 20
$ws.Range("F41").Value = 176
$ws.Range("F48").Value = 291
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2437
$ws.Range("F5").Value = 1652
$ws.Range("F15").Value = 322
$ws.Range("F16").Value = 2646
$ws.Range("F17").Value = 309
$ws.Range("F18").Value = 138
$ws.Range("F19").Value = 590
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2437
$ws.Range("F9").Value = 322
$ws.Range("F10").Value = 2646
$ws.Range("F11").Value = 309
$ws.Range("F13").Value = 7719
$ws.Range("F16").Value = 138
$ws.Range("F21").Value = 590
$ws.Range("F22").Value = 590
$ws.Range("F26").Value = 283
$ws.Range("F27").Value = 780
$ws.Range("F30").Value = 20
$ws.Range("F31").Value = 1344
$ws.Range("F35").Value = 531
$ws.Range("F39").Value = 254
$ws.Range("F50").Value = 291
